$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The date strings in column A (rows 3-21) use DD/MM/YYYY format values like
# "01/08/2022". Some of those, when assigned through .Value, would be
# auto-interpreted by Excel as real dates (since day <= 12 could also be a
# valid month). Force the range to Text format first so the dashed date
# strings are kept as plain text, matching the original inline-string cells.
$dateRange = $ws.Range("A3:A21")
$dateRange.NumberFormat = "@"

# Update date strings in column A (rows 3-21): replace "/" with "-"
for ($r = 3; $r -le 21; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $val = $cell.Value2
    if ($val -ne $null) {
        $cell.Value = $val -replace '/', '-'
    }
}

# Restore the default (Normal) style so no stray number-format style is left
# attached to these cells, matching the unchanged styling in the original file.
$dateRange.Style = "Normal"

# Update attendance counts for row 3: D3 and G3 go from 0 to 1
$ws.Range("D3").Value = 1
$ws.Range("G3").Value = 1
